$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J1 header cell - copy format from I1 (bold header style) then set text
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("J1").Value = "Code"

# Fill J column with "Dong_" + base filename values, in the original authoring order
# (row 29 was corrected last by the original author, producing the shared-string order in the diff)
$ws.Range("J4").Value = "Dong_Hmnioides"
$ws.Range("J5").Value = "Dong_Tlacunosa"
$ws.Range("J6").Value = "Dong_Apinguis"
$ws.Range("J7").Value = "Dong_Asharpii"
$ws.Range("J8").Value = "Dong_Bbarbata"
$ws.Range("J9").Value = "Dong_Bjaponica"
$ws.Range("J10").Value = "Dong_Btrilobata"
$ws.Range("J11").Value = "Dong_Cfissa"
$ws.Range("J12").Value = "Dong_Fcristula"
$ws.Range("J13").Value = "Dong_Forientalis"
$ws.Range("J14").Value = "Dong_Frullania"
$ws.Range("J15").Value = "Dong_Hdicranus"
$ws.Range("J16").Value = "Dong_Hramosus"
$ws.Range("J17").Value = "Dong_Hzollingeri"
$ws.Range("J18").Value = "Dong_Lsandvicensis"
$ws.Range("J19").Value = "Dong_Ltrichodes"
$ws.Range("J20").Value = "Dong_Mcrispata"
$ws.Range("J21").Value = "Dong_Malternifolia"
$ws.Range("J22").Value = "Dong_Mcrassipilis"
$ws.Range("J23").Value = "Dong_Mleptoneura"
$ws.Range("J24").Value = "Dong_Mnuda"
$ws.Range("J25").Value = "Dong_Ncurvifolia"
$ws.Range("J26").Value = "Dong_Ogrosseverrucosum"
$ws.Range("J27").Value = "Dong_Oprostratum"
$ws.Range("J28").Value = "Dong_Plyellii"
$ws.Range("J30").Value = "Dong_Pendiviifolia"
$ws.Range("J31").Value = "Dong_Pasplenioides"
$ws.Range("J32").Value = "Dong_Psubtropica"
$ws.Range("J33").Value = "Dong_Ppurpurea"
$ws.Range("J34").Value = "Dong_Phirtellus"
$ws.Range("J35").Value = "Dong_Pnavicularis"
$ws.Range("J36").Value = "Dong_Ppinnata"
$ws.Range("J37").Value = "Dong_Pplumosa"
$ws.Range("J38").Value = "Dong_Ppulcherrimum"
$ws.Range("J39").Value = "Dong_Pstriatus"
$ws.Range("J40").Value = "Dong_Rjaponica"
$ws.Range("J41").Value = "Dong_Rlindenbergia"
$ws.Range("J42").Value = "Dong_Rlatifrons"
$ws.Range("J43").Value = "Dong_Snemorosa"
$ws.Range("J44").Value = "Dong_Sornithopodioides"
$ws.Range("J45").Value = "Dong_Schistochila"
$ws.Range("J46").Value = "Dong_Ttomentella"
$ws.Range("J47").Value = "Dong_Awallichiana"
$ws.Range("J48").Value = "Dong_Bpusilla"
$ws.Range("J49").Value = "Dong_Cconicum"
$ws.Range("J50").Value = "Dong_Dhirsuta"
$ws.Range("J51").Value = "Dong_Lcruciata"
$ws.Range("J52").Value = "Dong_Mpaleacea"
$ws.Range("J53").Value = "Dong_Mpolymorpha"
$ws.Range("J54").Value = "Dong_Mtenerum"
$ws.Range("J55").Value = "Dong_Rberychiana"
$ws.Range("J56").Value = "Dong_Rcavernosa"
$ws.Range("J57").Value = "Dong_Stexanus"
$ws.Range("J58").Value = "Dong_Wdenudata"
$ws.Range("J29").Value = "Dong_Pepiphylla"

# Column J width
$ws.Columns.Item(10).ColumnWidth = 21.5

# Restore view state: scroll position and active selection
$ws.Range("A13").Select() | Out-Null
$ws.Range("J30").Select() | Out-Null
